$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The combined "file list" string that replaces the single "west.xlsx" value
# in F49, and is now also stamped (along with the source citation already
# used elsewhere in column E) across rows 50-93.
$combined = '"midwest.xlsx""northeast.xlsx""west.xlsx""south.xlsx"'
$source   = "Source: Consumer Expenditure Survey, U.S. Bureau of Labor Statistics, September, 2018"

# Row 49: F49 switches from the single "west.xlsx" citation to the combined
# four-file citation (new shared string).
$ws.Range("F49").Value = $combined

# Rows 50-93: stamp the same source (col E) / combined file list (col F)
# that the rows above them already carry.
for ($r = 50; $r -le 93; $r++) {
    $ws.Cells.Item($r, 5).Value = $source
    $ws.Cells.Item($r, 6).Value = $combined
}

# Column F needs to widen to fit the much longer combined string.
$ws.Columns("F").ColumnWidth = 44.71

# Restore the view state recorded in the saved file: scrolled so row 25 is
# at the top, with G49 as the active selection.
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("G49").Select() | Out-Null
